$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out the values that were mistakenly entered on Joshua's row (row 4)
$ws.Range("B4:H4").ClearContents()

# Correctly enter the availability on Adam's row (row 10)
$ws.Range("B10:F10").Value = "5pm-MN"
$ws.Range("G10:H10").Value = "8am-MN"

# Update the active selection to reflect the corrected cells
$ws.Activate()
$ws.Range("B10:H10").Select()
